# Generate Report for Handoff
# - Removes the stale "f0d6173c-..." row (row 3) from every sheet.
# - Updates the "ae18194b-..." row (row 2) to reflect the new handoff status/timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Materialize the hyperlink collection before deleting individual items,
# otherwise Delete() on a live COM enumeration removes everything.
$links1 = @($ws1.Hyperlinks)
for ($i = $links1.Count - 1; $i -ge 1; $i--) {
    $links1[$i].Delete()
}

# Drop the row for f0d6173c-... (row 3).
$ws1.Rows(3).Delete()

# Refresh status/date for the remaining ae18194b-... row.
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-07-17 06:07:31"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$links2 = @($ws2.Hyperlinks)
for ($i = $links2.Count - 1; $i -ge 5; $i--) {
    $links2[$i].Delete()
}

$ws2.Rows(3).Delete()

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-17 06:07:24"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$links3 = @($ws3.Hyperlinks)
for ($i = $links3.Count - 1; $i -ge 5; $i--) {
    $links3[$i].Delete()
}

$ws3.Rows(3).Delete()

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-17 06:07:31"
